# Commit: "updated slide templates + paths"
#
# Data change (the real content edit buried in the diff, once you account
# for shared-string re-indexing): in rows 38-52 of Sheet1, column E ("Stage")
# held the text "M" in every row. That literal text value is replaced with a
# genuine number: 4 for the Level-4/earlier rows (38-44) and 5 for the
# remaining rows (45-52). Removing the now-unused "M" shared string is just
# a side effect of that edit, which is why every shared-string index at/after
# the old "M" slot (and the cells that referenced them, e.g. column A/F)
# shifts down by one in the saved file - no other cell content actually
# changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E38:E44").Value = 4
$ws.Range("E45:E52").Value = 5
